# The deck ships with two themes:
#   ppt/theme/theme1.xml -> "Integral" (Red Violet clrScheme), used by the slide master
#   ppt/theme/theme2.xml -> "Office Theme" (Office clrScheme), used by the notes master
# The target edit swaps the two themes' contents: the slide master's theme
# becomes the "Office Theme" palette, and the notes master's theme becomes
# the "Integral" palette.
#
# RGB() isn't available in this host, so colours are passed as packed
# 0xBBGGRR integers (the same value PowerPoint's .RGB property reads/writes).

# ---- Slide master theme (ppt/theme/theme1.xml): Integral -> Office Theme ----
$smScheme = $ppt.ActivePresentation.SlideMaster.Theme.ThemeColorScheme
$smScheme.Item(1).RGB  = 0         # dk1      000000
$smScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$smScheme.Item(3).RGB  = 6968388   # dk2      44546A
$smScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$smScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$smScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$smScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$smScheme.Item(8).RGB  = 49407     # accent4  FFC000
$smScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$smScheme.Item(10).RGB = 4697456   # accent6  70AD47
$smScheme.Item(11).RGB = 12673797  # hlink    0563C1
$smScheme.Item(12).RGB = 7491477   # folHlink 954F72

# ---- Notes master theme (ppt/theme/theme2.xml): Office Theme -> Integral ----
$nmScheme = $ppt.ActivePresentation.NotesMaster.Theme.ThemeColorScheme
$nmScheme.Item(1).RGB  = 0         # dk1      000000
$nmScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$nmScheme.Item(3).RGB  = 5326149   # dk2      454551
$nmScheme.Item(4).RGB  = 14473688  # lt2      D8D9DC
$nmScheme.Item(5).RGB  = 9514467   # accent1  E32D91
$nmScheme.Item(6).RGB  = 13381832  # accent2  C830CC
$nmScheme.Item(7).RGB  = 14460494  # accent3  4EA6DC
$nmScheme.Item(8).RGB  = 15168839  # accent4  4775E7
$nmScheme.Item(9).RGB  = 14774665  # accent5  8971E1
$nmScheme.Item(10).RGB = 7555029   # accent6  D54773
$nmScheme.Item(11).RGB = 2465643   # hlink    6B9F25
$nmScheme.Item(12).RGB = 9211020   # folHlink 8C8C8C
